$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update case-count figures that changed for today's refresh ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1453381
$ws.Range("C4").Value = 23033
$ws.Range("D4").Value = 316244
$ws.Range("E4").Value = 1050367
$ws.Range("F4").Value = 16239
$ws.Range("G4").Value = 1573
$ws.Range("H4").Value = 86770

# Alemania (row 11)
$ws.Range("B11").Value = 174975
$ws.Range("C11").Value = 877
$ws.Range("E11").Value = 16747
$ws.Range("G11").Value = 67
$ws.Range("H11").Value = 7928

# Guinea (row 76)
$ws.Range("B76").Value = 2473
$ws.Range("C76").Value = 99
$ws.Range("D76").Value = 895
$ws.Range("E76").Value = 1563
$ws.Range("F76").Value = 18
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 15

# Guinea-Bisau (currently row 110) - rises above Tayikistan/Chipre/Albania/Libano/Niger
$ws.Range("B110").Value = 913
$ws.Range("C110").Value = 77
$ws.Range("D110").Value = 26
$ws.Range("E110").Value = 884
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 3

# Niger (currently row 109)
$ws.Range("B109").Value = 876
$ws.Range("C109").Value = 16
$ws.Range("D109").Value = 677
$ws.Range("E109").Value = 149
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 50

# Guinea Ecuatorial (row 123)
$ws.Range("B123").Value = 583
$ws.Range("C123").Value = 61
$ws.Range("E123").Value = 563
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 7

# Islas Caimanes (row 168)
$ws.Range("B168").Value = 93
$ws.Range("C168").Value = 7
$ws.Range("E168").Value = 38

# --- Timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Mayo de 2020 a las 00:05"

# --- Re-sort the country table by total cases (column B), descending ---
# (stable sort: ties keep their previous relative order)
$rng = $ws.Range("A4:H222")
$rng.Sort($ws.Range("B4:B222"), 2)

# --- San Bartolome / Sahara Occidental are tied on every figure; Sahara
#     Occidental now sorts ahead of San Bartolome, so swap those two rows ---
$a = $ws.Range("A215").Value2
$b = $ws.Range("B215").Value2
$c = $ws.Range("C215").Value2
$d = $ws.Range("D215").Value2
$e = $ws.Range("E215").Value2
$f = $ws.Range("F215").Value2
$g = $ws.Range("G215").Value2
$h = $ws.Range("H215").Value2

$ws.Range("A215").Value = $ws.Range("A216").Value2
$ws.Range("B215").Value = $ws.Range("B216").Value2
$ws.Range("C215").Value = $ws.Range("C216").Value2
$ws.Range("D215").Value = $ws.Range("D216").Value2
$ws.Range("E215").Value = $ws.Range("E216").Value2
$ws.Range("F215").Value = $ws.Range("F216").Value2
$ws.Range("G215").Value = $ws.Range("G216").Value2
$ws.Range("H215").Value = $ws.Range("H216").Value2

$ws.Range("A216").Value = $a
$ws.Range("B216").Value = $b
$ws.Range("C216").Value = $c
$ws.Range("D216").Value = $d
$ws.Range("E216").Value = $e
$ws.Range("F216").Value = $f
$ws.Range("G216").Value = $g
$ws.Range("H216").Value = $h
